# "changed batch size for aligned data"
# Adds a new "MSE:" column header (column D) to each of the six
# comparison tables on the Results sheet, and updates the view's
# scroll/selection state to reflect the newly active cell D36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("D1").Value  = "MSE:"
$ws.Range("D6").Value  = "MSE:"
$ws.Range("D14").Value = "MSE:"
$ws.Range("D22").Value = "MSE:"
$ws.Range("D29").Value = "MSE:"
$ws.Range("D36").Value = "MSE:"

# Scroll the view so row 22 is at the top and select D36, matching
# the saved view state in the edited workbook.
$excel.Goto($ws.Range("A22"), $true)
$ws.Range("D36").Select()
